$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# Update patient record id (A2) and ID card / insurance number (E2)
$ws.Range("A2").Value = 3019
$ws.Range("E2").Value = 46200608019

# Move the active selection to E2, matching the saved cursor position
$ws.Range("E2").Select() | Out-Null
